$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CapitalCall")

# Create the new "Exchange Rates" sheet right after "CapitalCall"
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Exchange Rates"

# Move the exchange-rate columns (J:M) from CapitalCall into the new sheet (A:D)
$src = $ws1.Range("J1:M4")
$src.Copy($newSheet.Range("A1"))

# Remove the now-relocated columns from the original sheet
$ws1.Range("J1:M4").Clear()

# Restore the selections seen in the target workbook
$newSheet.Range("D31").Select()
$ws1.Select()
$ws1.Range("C26").Select()
